$d = $word.ActiveDocument

# --- 1. Mark the two "noProof" drawings (paragraphs holding the floating
#        pictures inserted after the combobox / registration screenshots)
#        as not to be proofed, matching the <w:noProof/> added in rPr.
$pDraw1 = $d.Paragraphs.Item(76)
$pDraw1.Range.NoProofing = 1

$pDraw2 = $d.Paragraphs.Item(86)
$pDraw2.Range.NoProofing = 1

# --- 2. Locate the "Foto del estudiante" paragraph and append a colon as
#        its own run, then add a new paragraph with the student's name.
$target = $null
$targetIndex = 0
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "Foto del estudiante*") {
        $target = $p
        $targetIndex = $idx
    }
}

$insertPoint = $d.Range($target.Range.End - 1, $target.Range.End - 1)
$insertPoint.InsertAfter(":")
# Toggling formatting off forces the simulator to close the previous run
# and start a fresh one, matching the separate <w:r> for ":" in the diff.
$insertPoint.Font.Bold = 1
$insertPoint.Font.Bold = 0

$target.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.InsertAfter("Luis Carlos Salazar Quesquén – U22238714")
